$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.116.48'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '3.055.41'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.86'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.19'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").Value = '3.055.21'
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("E10").Value = '  -2.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.86'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("E12").Value = '  -2.75%  '
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("E14").Value = '  -3.28%  '
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("D16").Value = '3.560.04'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("D18").Value = '63.141.84'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '3.054.39'
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.16'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.71'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("E27").Value = '  +4.18%  '
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.37'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.64'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.49%  '
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("D36").Value = '0.0₃0821'
$ws.Range("E36").Value = '  -3.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.27'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.24'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '435.03'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.290'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("E44").Value = '  +3.00%  '
$ws.Range("E45").Value = '  -0.52%  '
$ws.Range("D46").Value = '2.828.42'
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.24'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.51'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.62%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.18'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("E51").Value = '  -1.59%  '
